$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 30

    $ws.Cells.Item($row, 1).Value = 29

    # Date-looking text must be forced to stay as literal text (not an
    # auto-converted date serial), matching the rest of column B.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "12:38:01"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.66

    # Trade is still open, so Exit Price is blank (present empty cell).
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.6550058009231
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Trade is still open, so Exit Reason is blank (present empty cell).
    $ws.Cells.Item($row, 16).Style = "Normal"

    $ws.Cells.Item($row, 17).Value = 0
}
